# The paragraph with the transcription tag currently reads the id value
# "<id>p144r_1</id>" split across three separately-formatted runs:
#   1) "<id>"      -> Courier New, color 7f6000, sz/szCs 18
#   2) "p144r_1"   -> default font, color 000000
#   3) "</id>"     -> Courier New, color 7f6000, sz/szCs 18
#
# The edit collapses those three runs into a single run containing the
# whole string "<id>p144r_1</id>", taking on the Courier-New/olive
# formatting that both the leading "<id>" run and trailing "</id>" run
# already share.
#
# Doing a plain Find/Replace over the full "<id>p144r_1</id>" span lets
# Word re-flow the replacement text using the formatting already present
# at the (identically formatted) boundary of the match, which naturally
# fuses the three runs into one - exactly mirroring the target diff.

$d = $word.ActiveDocument
$rng = $d.Content

$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()

$found = $rng.Find.Execute(
    "<id>p144r_1</id>",  # Find
    $true,               # MatchCase
    $false,              # MatchWholeWord
    $false,              # MatchWildcards
    $false,              # MatchSoundsLike
    $false,              # MatchAllWordForms
    $true,               # Forward
    1,                   # Wrap           (wdFindContinue)
    $false,              # Format
    "<id>p144r_1</id>",  # ReplaceWith
    1                    # Replace        (wdReplaceOne)
)

if (-not $found) {
    Write-Host "WARNING: target text '<id>p144r_1</id>' was not found"
}
